$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "促甲状腺素"
$ws.Range("A3").Value = "血清甲状腺素"
$ws.Range("A4").Value = "血清游离四碘甲状腺原氨酸"
$ws.Range("C4").Value = "pg"
$ws.Range("A5").Value = "血清三碘甲状腺原氨酸"
$ws.Range("A6").Value = "血清游离三碘甲状腺原氨酸"
$ws.Range("C6").Value = "pg"
